$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.1
$ws.Range("G2").Value = 2.6
$ws.Range("H2").Value = 2.66
$ws.Range("J2").Value = 3.35
$ws.Range("K2").Value = 4.9
$ws.Range("L2").Value = 1.21
$ws.Range("N2").Value = 2.08
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.54
$ws.Range("R2").Value = 1.41
$ws.Range("S2").Value = 2.28
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.33
$ws.Range("W2").Value = 1.71

# Row 3
$ws.Range("K3").Value = 5.1
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 1.16
$ws.Range("P3").Value = 2.74
$ws.Range("Q3").Value = 1.49
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = 2.22
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 2.16
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 2.6
$ws.Range("X3").Value = 29
$ws.Range("Y3").Value = 32
$ws.Range("Z3").Value = 60
$ws.Range("AA3").Value = 160
$ws.Range("AB3").Value = 15
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 27
$ws.Range("AE3").Value = 60
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 10.5
$ws.Range("AH3").Value = 18
$ws.Range("AI3").Value = 55
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 15
$ws.Range("AL3").Value = 25
$ws.Range("AM3").Value = 70
$ws.Range("AN3").Value = 6.2
$ws.Range("AO3").Value = 50

# Row 4
$ws.Range("F4").Value = 2.74
$ws.Range("H4").Value = 2.44
$ws.Range("I4").Value = 2.8
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 2.18
$ws.Range("O4").Value = 1.22
$ws.Range("Q4").Value = 1.69
$ws.Range("R4").Value = 1.41
$ws.Range("S4").Value = 2.46
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.55
$ws.Range("W4").Value = 1.48
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 20
$ws.Range("Z4").Value = 28
$ws.Range("AA4").Value = 50
$ws.Range("AB4").Value = 21
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 18
$ws.Range("AE4").Value = 38
$ws.Range("AF4").Value = 30
$ws.Range("AG4").Value = 19
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 50
$ws.Range("AJ4").Value = 65
$ws.Range("AK4").Value = 44
$ws.Range("AL4").Value = 50
$ws.Range("AM4").Value = 100
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("G5").Value = 3.65
$ws.Range("H5").Value = 2.04
$ws.Range("K5").Value = 4.6
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 2.66
$ws.Range("O5").Value = 1.15
$ws.Range("R5").Value = 1.59
$ws.Range("S5").Value = 2.1
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.87
$ws.Range("W5").Value = 1.37
$ws.Range("X5").Value = 36
$ws.Range("Y5").Value = 21
$ws.Range("Z5").Value = 24
$ws.Range("AA5").Value = 36
$ws.Range("AB5").Value = 30
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 16
$ws.Range("AE5").Value = 26
$ws.Range("AF5").Value = 42
$ws.Range("AG5").Value = 22
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 36
$ws.Range("AJ5").Value = 80
$ws.Range("AK5").Value = 48
$ws.Range("AL5").Value = 48
$ws.Range("AM5").Value = 75
$ws.Range("AN5").Value = 30
$ws.Range("AO5").Value = 12.5

# Row 6
$ws.Range("G6").Value = 7.6
$ws.Range("H6").Value = 1.49
$ws.Range("J6").Value = 4.8
$ws.Range("L6").Value = 1.28
$ws.Range("Q6").Value = 1.61
$ws.Range("T6").Value = 1.66
$ws.Range("V6").Value = 2.92
$ws.Range("W6").Value = 1.15
$ws.Range("X6").Value = 23
$ws.Range("Z6").Value = 11
$ws.Range("AC6").Value = 11.5
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 14.5
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 25
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 220
$ws.Range("AK6").Value = 95
$ws.Range("AL6").Value = 85
$ws.Range("AM6").Value = 110
$ws.Range("AN6").Value = 100
$ws.Range("AO6").Value = 6.8
